$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for the Price column (D) so numeric-looking
# strings like "1.001" or "30.444.81" are not auto-converted to numbers.

$dStyle2 = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.444.81"
$ws.Range("D2").Style = $dStyle2
$ws.Range("E2").Value = "  +1.02%  "

$dStyle3 = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.872.38"
$ws.Range("D3").Style = $dStyle3
$ws.Range("E3").Value = "  +0.62%  "

$dStyle4 = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = $dStyle4
$ws.Range("E4").Value = "  +0.25%  "

$dStyle5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.88"
$ws.Range("D5").Style = $dStyle5
$ws.Range("E5").Value = "  +2.21%  "

$dStyle6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = $dStyle6
$ws.Range("E6").Value = "  +0.28%  "

$dStyle7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4737"
$ws.Range("D7").Style = $dStyle7
$ws.Range("E7").Value = "  +1.13%  "

$dStyle8 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2908"
$ws.Range("D8").Style = $dStyle8
$ws.Range("E8").Value = "  +1.77%  "

$dStyle9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06502"
$ws.Range("D9").Style = $dStyle9
$ws.Range("E9").Value = "  +0.44%  "

$dStyle10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.95"
$ws.Range("D10").Style = $dStyle10
$ws.Range("E10").Value = "  +5.43%  "

$dStyle11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07718"
$ws.Range("D11").Style = $dStyle11
$ws.Range("E11").Value = "  +0.48%  "

$dStyle12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.63"
$ws.Range("D12").Style = $dStyle12
$ws.Range("E12").Value = "  +4.03%  "

$dStyle13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7398"
$ws.Range("D13").Style = $dStyle13
$ws.Range("E13").Value = "  +8.83%  "

$dStyle14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.882.55"
$ws.Range("D14").Style = $dStyle14
$ws.Range("E14").Value = "  +1.17%  "

$dStyle15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.122"
$ws.Range("D15").Style = $dStyle15
$ws.Range("E15").Value = "  +1.08%  "

$dStyle16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "273.33"
$ws.Range("D16").Style = $dStyle16
$ws.Range("E16").Value = "  +1.53%  "

$dStyle17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.425.55"
$ws.Range("D17").Style = $dStyle17
$ws.Range("E17").Value = "  +1.04%  "

$dStyle18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.40"
$ws.Range("D18").Style = $dStyle18
$ws.Range("E18").Value = "  +0.67%  "

$dStyle19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007535"
$ws.Range("D19").Style = $dStyle19
$ws.Range("E19").Value = "  +0.24%  "

$dStyle20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").Style = $dStyle20
$ws.Range("E20").Value = "  +0.31%  "

$dStyle21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.127.79"
$ws.Range("D21").Style = $dStyle21
$ws.Range("E21").Value = "  +1.74%  "

$dStyle22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = $dStyle22
$ws.Range("E22").Value = "  +0.27%  "

$dStyle23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.227"
$ws.Range("D23").Style = $dStyle23
$ws.Range("E23").Value = "  +1.32%  "

$dStyle24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.175"
$ws.Range("D24").Style = $dStyle24
$ws.Range("E24").Value = "  +1.01%  "

$dStyle25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.254"
$ws.Range("D25").Style = $dStyle25
$ws.Range("E25").Value = "  -0.42%  "

$dStyle26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.85"
$ws.Range("D26").Style = $dStyle26
$ws.Range("E26").Value = "  -1.22%  "

$dStyle27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.83"
$ws.Range("D27").Style = $dStyle27
$ws.Range("E27").Value = "  +0.32%  "

$dStyle28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.932"
$ws.Range("D28").Style = $dStyle28
$ws.Range("E28").Value = "  +2.59%  "

$dStyle29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1006"
$ws.Range("D29").Style = $dStyle29
$ws.Range("E29").Value = "  +2.50%  "

$dStyle30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.365"
$ws.Range("D30").Style = $dStyle30
$ws.Range("E30").Value = "  -0.30%  "

$dStyle31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.511"
$ws.Range("D31").Style = $dStyle31
$ws.Range("E31").Value = "  +0.74%  "

$dStyle32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.315"
$ws.Range("D32").Style = $dStyle32
$ws.Range("E32").Value = "  +2.49%  "

$dStyle33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.102"
$ws.Range("D33").Style = $dStyle33
$ws.Range("E33").Value = "  +2.85%  "

$dStyle34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04817"
$ws.Range("D34").Style = $dStyle34
$ws.Range("E34").Value = "  +3.25%  "

$dStyle35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.123"
$ws.Range("D35").Style = $dStyle35
$ws.Range("E35").Value = "  +1.32%  "

$dStyle36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6984"
$ws.Range("D36").Style = $dStyle36
$ws.Range("E36").Value = "  +1.90%  "

$dStyle37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("D37").Style = $dStyle37
$ws.Range("E37").Value = "  +0.17%  "

$dStyle38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.721"
$ws.Range("D38").Style = $dStyle38
$ws.Range("E38").Value = "  +0.62%  "

$dStyle39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01860"
$ws.Range("D39").Style = $dStyle39
$ws.Range("E39").Value = "  +1.88%  "

$dStyle40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.740"
$ws.Range("D40").Style = $dStyle40
$ws.Range("E40").Value = "  +0.87%  "

$dStyle41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.287"
$ws.Range("D41").Style = $dStyle41
$ws.Range("E41").Value = "  -0.59%  "

$dStyle42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.978"
$ws.Range("D42").Style = $dStyle42
$ws.Range("E42").Value = "  +4.92%  "

$dStyle43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.96"
$ws.Range("D43").Style = $dStyle43
$ws.Range("E43").Value = "  +0.91%  "

$dStyle44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4193"
$ws.Range("D44").Style = $dStyle44
$ws.Range("E44").Value = "  +3.65%  "

$dStyle45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = $dStyle45
$ws.Range("E45").Value = "  +0.31%  "

$dStyle46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.8368"
$ws.Range("D46").Style = $dStyle46
$ws.Range("E46").Value = "  +0.48%  "

$dStyle47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.64"
$ws.Range("D47").Style = $dStyle47
$ws.Range("E47").Value = "  +0.83%  "

$dStyle48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.290"
$ws.Range("D48").Style = $dStyle48
$ws.Range("E48").Value = "  +0.97%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$dStyle49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.53"
$ws.Range("D49").Style = $dStyle49
$ws.Range("E49").Value = "  +4.08%  "

$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$dStyle50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.019"
$ws.Range("D50").Style = $dStyle50
$ws.Range("E50").Value = "  +1.19%  "

$dStyle51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "919.32"
$ws.Range("D51").Style = $dStyle51
$ws.Range("E51").Value = "  -0.96%  "
